# This script applies the stock-report refresh edit described by the commit diff.
# It updates quantity (F) / value (G) figures for re-counted items, recomputes the
# affected company Sub Total (B) and Grand Total (B718/B719) rows, swaps a handful of
# duplicate-item row pairs whose order was corrected, and fixes one item-name casing.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F20").Value = 331
$ws.Range("G20").Value = 16957.13
$ws.Range("B34").Value = 60115.82
$ws.Range("F44").Value = 35
$ws.Range("G44").Value = 1236.2
$ws.Range("F51").Value = 151
$ws.Range("G51").Value = 14124.54
$ws.Range("F58").Value = 76
$ws.Range("G58").Value = 5922.68
$ws.Range("F65").Value = 11
$ws.Range("G65").Value = 368.06
$ws.Range("B66").Value = 210750.8
$ws.Range("F105").Value = 72
$ws.Range("G105").Value = 5844.24
$ws.Range("F106").Value = 186
$ws.Range("G106").Value = 20865.48
$ws.Range("F120").Value = 89
$ws.Range("G120").Value = 11994.53
$ws.Range("B123").Value = 75423.09
$ws.Range("B126").Value = 64196
$ws.Range("F126").Value = 1
$ws.Range("G126").Value = 32143.58
$ws.Range("B127").Value = 65258
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("F139").Value = 28
$ws.Range("G139").Value = 1907.64
$ws.Range("F141").Value = 55
$ws.Range("G141").Value = 2863.3
$ws.Range("B147").Value = 22173.51
$ws.Range("B161").Value = 64350
$ws.Range("E161").Value = 70.63
$ws.Range("F161").Value = 2
$ws.Range("G161").Value = 132.88
$ws.Range("B162").Value = 57756
$ws.Range("E162").Value = 79.37
$ws.Range("F162").Value = -100
$ws.Range("G162").Value = -6644
$ws.Range("F173").Value = 52
$ws.Range("G173").Value = 4087.72
$ws.Range("F186").Value = 22
$ws.Range("G186").Value = 952.16
$ws.Range("F192").Value = 26
$ws.Range("G192").Value = 1557.14
$ws.Range("B193").Value = 67434.86
$ws.Range("F203").Value = 4
$ws.Range("G203").Value = 33.8
$ws.Range("B204").Value = 33.8
$ws.Range("F206").Value = 73
$ws.Range("G206").Value = 4730.4
$ws.Range("B208").Value = 4776.89
$ws.Range("F210").Value = 131
$ws.Range("G210").Value = 7127.71
$ws.Range("F213").Value = 220
$ws.Range("G213").Value = 27869.6
$ws.Range("F214").Value = 101
$ws.Range("G214").Value = 10651.46
$ws.Range("B218").Value = 82225.2
$ws.Range("F222").Value = 949
$ws.Range("G222").Value = 17556.5
$ws.Range("F223").Value = 38
$ws.Range("G223").Value = 814.34
$ws.Range("B229").Value = 29683.42
$ws.Range("F262").Value = 77
$ws.Range("G262").Value = 6031.41
$ws.Range("F263").Value = 16
$ws.Range("G263").Value = 1659.2
$ws.Range("F264").Value = 77
$ws.Range("G264").Value = 2682.68
$ws.Range("F266").Value = 8
$ws.Range("G266").Value = 532.72
$ws.Range("F274").Value = 51
$ws.Range("G274").Value = 1776.84
$ws.Range("F285").Value = 11
$ws.Range("G285").Value = 1220.34
$ws.Range("B290").Value = 66194
$ws.Range("C290").Value = "HIM-Total Care Baby Pants Diapers-M-9s"
$ws.Range("F290").Value = 27
$ws.Range("G290").Value = 2313.36
$ws.Range("B291").Value = 64983
$ws.Range("C291").Value = "HIM-TOTAL CARE BABY PANTS DIAPERS-M-9S"
$ws.Range("F291").Value = 6
$ws.Range("G291").Value = 514.08
$ws.Range("B295").Value = 124714.42
$ws.Range("B308").Value = 55356
$ws.Range("E308").Value = 54.04
$ws.Range("F308").Value = -158
$ws.Range("G308").Value = -7527.12
$ws.Range("B309").Value = 63510
$ws.Range("E309").Value = 50.66
$ws.Range("F309").Value = 78
$ws.Range("G309").Value = 3715.92
$ws.Range("F324").Value = 54
$ws.Range("G324").Value = 9251.82
$ws.Range("F325").Value = 39
$ws.Range("G325").Value = 5895.63
$ws.Range("B328").Value = -3369.14
$ws.Range("F349").Value = 152
$ws.Range("G349").Value = 11343.76
$ws.Range("B356").Value = 79364.51
$ws.Range("F361").Value = 242
$ws.Range("G361").Value = 34022.78
$ws.Range("B363").Value = 78532.79
$ws.Range("F368").Value = 60
$ws.Range("G368").Value = 1921.2
$ws.Range("F371").Value = 66
$ws.Range("G371").Value = 9917.82
$ws.Range("B372").Value = 64043.79
$ws.Range("B381").Value = 58047
$ws.Range("D381").Value = 105.54
$ws.Range("E381").Value = 126.1
$ws.Range("F381").Value = 32
$ws.Range("G381").Value = 3377.28
$ws.Range("B382").Value = 47097
$ws.Range("D382").Value = 112.28
$ws.Range("E382").Value = 134.16
$ws.Range("F382").Value = 15
$ws.Range("G382").Value = 1684.2
$ws.Range("F387").Value = 438
$ws.Range("G387").Value = 42310.8
$ws.Range("B389").Value = 59252.26
$ws.Range("F396").Value = 135
$ws.Range("G396").Value = 3439.8
$ws.Range("B417").Value = 174201.5
$ws.Range("F433").Value = 138
$ws.Range("G433").Value = 1330.32
$ws.Range("B438").Value = 25964.74
$ws.Range("F453").Value = 40
$ws.Range("G453").Value = 5812.8
$ws.Range("F456").Value = 149
$ws.Range("G456").Value = 40210.63
$ws.Range("B458").Value = 100076.9
$ws.Range("F478").Value = 11
$ws.Range("G478").Value = 2439.14
$ws.Range("B479").Value = 53319
$ws.Range("E479").Value = 310.64
$ws.Range("F479").Value = -6
$ws.Range("G479").Value = -1643.52
$ws.Range("B480").Value = 64810
$ws.Range("E480").Value = 291.22
$ws.Range("F480").Value = 0
$ws.Range("G480").Value = 0
$ws.Range("B482").Value = 2756.9
$ws.Range("F511").Value = 250
$ws.Range("G511").Value = 24967.5
$ws.Range("F512").Value = 24
$ws.Range("G512").Value = 2845.92
$ws.Range("F519").Value = 420
$ws.Range("G519").Value = 23049.6
$ws.Range("F520").Value = 33
$ws.Range("G520").Value = 904.2
$ws.Range("F521").Value = 97
$ws.Range("G521").Value = 2657.8
$ws.Range("F522").Value = 88
$ws.Range("G522").Value = 2344.32
$ws.Range("B525").Value = 129676.43
$ws.Range("F527").Value = 55
$ws.Range("G527").Value = 1821.05
$ws.Range("F528").Value = 291
$ws.Range("G528").Value = 4615.26
$ws.Range("F531").Value = 221
$ws.Range("G531").Value = 7317.31
$ws.Range("B535").Value = 25180.57
$ws.Range("F544").Value = 43
$ws.Range("G544").Value = 2661.7
$ws.Range("B556").Value = 51256.78
$ws.Range("F558").Value = 206
$ws.Range("G558").Value = 25101.1
$ws.Range("B561").Value = 29591.3
$ws.Range("F616").Value = 4
$ws.Range("G616").Value = 571.08
$ws.Range("F617").Value = 25
$ws.Range("G617").Value = 1203
$ws.Range("F620").Value = 370
$ws.Range("G620").Value = 29078.3
$ws.Range("F621").Value = 20
$ws.Range("G621").Value = 7556.2
$ws.Range("F625").Value = 328
$ws.Range("G625").Value = 12080.24
$ws.Range("B628").Value = 213277.83
$ws.Range("F659").Value = 40
$ws.Range("G659").Value = 2141.6
$ws.Range("B668").Value = 12932.19
$ws.Range("B718").Value = 2816432.47
$ws.Range("B719").Value = 2816432.47
